$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = -7.887
$ws.Range("B12").Value = 6.112
$ws.Range("D12").Value = -8.311
$ws.Range("D14").Value = -8.263
$ws.Range("D22").Value = -8.106
